$wb = $excel.ActiveWorkbook

# Generate Report for Handoff:
# Update the handoff timestamps for the b413f230-3044-4ce5-91a8-1e66ce69384d.md
# file (row 7 on every sheet) to reflect the freshly generated xliff hand-off.

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("G7").Value = "2016-08-23 02:41:19"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("H7").Value = "2016-08-23 02:41:14"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("H7").Value = "2016-08-23 02:41:19"
